# Add database connection for component.equipment_penetration_rates
# Remove the "[m/h]" unit suffix from the technique labels on the "penet"
# sheet so the values line up with the new database connection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("penet")

$ws.Range("A2").Value = "Drilling rig"
$ws.Range("A3").Value = "Hammer"
$ws.Range("A4").Value = "Vibro driver"
$ws.Range("A5").Value = "ROV with suction pump"
$ws.Range("A6").Value = "ROV with jetting"

$ws.Range("A6").Select()
